$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1150.2285  # H15 was 1474.7941
$ws.Cells.Item(15, 9).Value = 1150.2285  # I15 was 1474.7941
$ws.Cells.Item(15, 11).Value = 3450.6855  # K15 was 4424.3823
$ws.Cells.Item(15, 13).Value = -3281.6855  # M15 was -4255.3823
$ws.Cells.Item(33, 8).Value = 958568.3  # H33 was 1015450.56
$ws.Cells.Item(33, 9).Value = 1014690  # I33 was 1149944
$ws.Cells.Item(33, 10).Value = 4500  # J33 was 6750
$ws.Cells.Item(33, 11).Value = 1014690  # K33 was 1149944
$ws.Cells.Item(33, 12).Value = 4500  # L33 was 6750
$ws.Cells.Item(33, 13).Value = -1014461  # M33 was -1149715
$ws.Cells.Item(33, 14).Value = -4958  # N33 was -7208
$ws.Cells.Item(64, 8).Value = 7533.1665  # H64 was 6374.875
$ws.Cells.Item(64, 9).Value = 4099.5  # I64 was 3666.3333
$ws.Cells.Item(64, 10).Value = 9250  # J64 was 8000
$ws.Cells.Item(64, 11).Value = 4099.5  # K64 was 3666.3333
$ws.Cells.Item(64, 12).Value = 9250  # L64 was 8000
$ws.Cells.Item(64, 13).Value = -3851.5  # M64 was -3418.3333
$ws.Cells.Item(64, 14).Value = -9746  # N64 was -8496
$ws.Cells.Item(67, 8).Value = 7533.1665  # H67 was 6374.875
$ws.Cells.Item(67, 9).Value = 4099.5  # I67 was 3666.3333
$ws.Cells.Item(67, 10).Value = 9250  # J67 was 8000
$ws.Cells.Item(67, 11).Value = 4099.5  # K67 was 3666.3333
$ws.Cells.Item(67, 12).Value = 9250  # L67 was 8000
$ws.Cells.Item(67, 13).Value = -3241.5  # M67 was -2808.3333
$ws.Cells.Item(67, 14).Value = -10966  # N67 was -9716
$ws.Cells.Item(87, 8).Value = 33354  # H87 was 33340.5
$ws.Cells.Item(87, 10).Value = 33354  # J87 was 33340.5
$ws.Cells.Item(87, 12).Value = 33354  # L87 was 33340.5
$ws.Cells.Item(87, 14).Value = -35850  # N87 was -35836.5
$ws.Cells.Item(90, 8).Value = 33354  # H90 was 33340.5
$ws.Cells.Item(90, 10).Value = 33354  # J90 was 33340.5
$ws.Cells.Item(90, 12).Value = 100062  # L90 was 100021.5
$ws.Cells.Item(90, 14).Value = -112542  # N90 was -112501.5
$ws.Cells.Item(129, 8).Value = 873.4375  # H129 was 937.5
$ws.Cells.Item(129, 9).Value = 741.0714  # I129 was 793.75
$ws.Cells.Item(129, 11).Value = 2223.2142  # K129 was 2381.25
$ws.Cells.Item(129, 13).Value = 2776.7858  # M129 was 2618.75
$ws.Cells.Item(132, 8).Value = 3397.8293  # H132 was 3160
$ws.Cells.Item(132, 9).Value = 1267.2424  # I132 was 1222.8611
$ws.Cells.Item(132, 10).Value = 12186.5  # J132 was 10908.556
$ws.Cells.Item(132, 11).Value = 3801.7272  # K132 was 3668.5833
$ws.Cells.Item(132, 12).Value = 36559.5  # L132 was 32725.668
$ws.Cells.Item(132, 13).Value = -1271.7272  # M132 was -1138.5833
$ws.Cells.Item(132, 14).Value = -41619.5  # N132 was -37785.66800000001
$ws.Cells.Item(137, 8).Value = 36294.137  # H137 was 38880.742
$ws.Cells.Item(137, 9).Value = 1681.1818  # I137 was 1716
$ws.Cells.Item(137, 10).Value = 145077.72  # J137 was 168957.33
$ws.Cells.Item(137, 11).Value = 5043.5454  # K137 was 5148
$ws.Cells.Item(137, 12).Value = 435233.16  # L137 was 506871.99
$ws.Cells.Item(137, 13).Value = -2493.5454  # M137 was -2598
$ws.Cells.Item(137, 14).Value = -440333.16  # N137 was -511971.99
$ws.Cells.Item(138, 8).Value = 2823.3572  # H138 was 2484.4546
$ws.Cells.Item(138, 9).Value = 1693  # I138 was 1770
$ws.Cells.Item(138, 10).Value = 5649.25  # J138 was 5699.5
$ws.Cells.Item(138, 11).Value = 5079  # K138 was 5310
$ws.Cells.Item(138, 12).Value = 16947.75  # L138 was 17098.5
$ws.Cells.Item(138, 13).Value = 61  # M138 was -170
$ws.Cells.Item(138, 14).Value = -27227.75  # N138 was -27378.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 133758.4  # H4 was 143144.86
$ws.Cells.Item(4, 9).Value = 143304  # I4 was 166983.58
$ws.Cells.Item(4, 10).Value = 120  # J4 was 112.5
$ws.Cells.Item(4, 11).Value = 143304  # K4 was 166983.58
$ws.Cells.Item(4, 12).Value = 120  # L4 was 112.5
$ws.Cells.Item(4, 13).Value = -143188  # M4 was -166867.58
$ws.Cells.Item(4, 14).Value = -352  # N4 was -344.5
$ws.Cells.Item(5, 8).Value = 133.45454  # H5 was 126.083336
$ws.Cells.Item(5, 9).Value = 143.28572  # I5 was 118.111115
$ws.Cells.Item(5, 10).Value = 116.25  # J5 was 150
$ws.Cells.Item(5, 11).Value = 143.28572  # K5 was 118.111115
$ws.Cells.Item(5, 12).Value = 116.25  # L5 was 150
$ws.Cells.Item(5, 13).Value = -31.28572  # M5 was -6.111114999999998
$ws.Cells.Item(5, 14).Value = -340.25  # N5 was -374
$ws.Cells.Item(6, 8).Value = 10000  # H6 was 0
$ws.Cells.Item(6, 9).Value = 10000  # I6 was 0
$ws.Cells.Item(6, 11).Value = 10000  # K6 was 0
$ws.Cells.Item(6, 13).Value = -9827  # M6 was None
$ws.Cells.Item(32, 8).Value = 38308.242  # H32 was 39042.848
$ws.Cells.Item(32, 9).Value = 23875.045  # I32 was 24428.559
$ws.Cells.Item(32, 10).Value = 101814.3  # J32 was 101884.3
$ws.Cells.Item(32, 11).Value = 23875.045  # K32 was 24428.559
$ws.Cells.Item(32, 12).Value = 101814.3  # L32 was 101884.3
$ws.Cells.Item(32, 13).Value = -23588.045  # M32 was -24141.559
$ws.Cells.Item(32, 14).Value = -102388.3  # N32 was -102458.3
$ws.Cells.Item(113, 8).Value = 99499.5  # H113 was 99999
$ws.Cells.Item(113, 10).Value = 99499.5  # J113 was 99999
$ws.Cells.Item(113, 12).Value = 99499.5  # L113 was 99999
$ws.Cells.Item(113, 14).Value = -108177.5  # N113 was -108677
$ws.Cells.Item(122, 8).Value = 22408.9  # H122 was 16411.715
$ws.Cells.Item(122, 9).Value = 30941.285  # I122 was 22088.4
$ws.Cells.Item(122, 10).Value = 2500  # J122 was 2220
$ws.Cells.Item(122, 11).Value = 92823.855  # K122 was 66265.20000000001
$ws.Cells.Item(122, 12).Value = 7500  # L122 was 6660
$ws.Cells.Item(122, 13).Value = -90373.855  # M122 was -63815.20000000001
$ws.Cells.Item(122, 14).Value = -12400  # N122 was -11560

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 133.45454  # H4 was 126.083336
$ws.Cells.Item(4, 9).Value = 143.28572  # I4 was 118.111115
$ws.Cells.Item(4, 10).Value = 116.25  # J4 was 150
$ws.Cells.Item(4, 11).Value = 143.28572  # K4 was 118.111115
$ws.Cells.Item(4, 12).Value = 116.25  # L4 was 150
$ws.Cells.Item(4, 13).Value = -28.28572  # M4 was -3.111114999999998
$ws.Cells.Item(4, 14).Value = -346.25  # N4 was -380
$ws.Cells.Item(86, 8).Value = 1696.5  # H86 was 1755.2
$ws.Cells.Item(86, 9).Value = 1429.125  # I86 was 1474.7858
$ws.Cells.Item(86, 11).Value = 1429.125  # K86 was 1474.7858
$ws.Cells.Item(86, 13).Value = -306.125  # M86 was -351.7858000000001
$ws.Cells.Item(89, 8).Value = 1696.5  # H89 was 1755.2
$ws.Cells.Item(89, 9).Value = 1429.125  # I89 was 1474.7858
$ws.Cells.Item(89, 11).Value = 7145.625  # K89 was 7373.929
$ws.Cells.Item(89, 13).Value = -1529.625  # M89 was -1757.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2008.5714  # H58 was 1443.5
$ws.Cells.Item(58, 9).Value = 1510  # I58 was 1443.5
$ws.Cells.Item(58, 10).Value = 5000  # J58 was 0
$ws.Cells.Item(58, 11).Value = 1510  # K58 was 1443.5
$ws.Cells.Item(58, 12).Value = 5000  # L58 was 0
$ws.Cells.Item(58, 13).Value = -1307  # M58 was -1240.5
$ws.Cells.Item(58, 14).Value = -5406  # N58 was None
$ws.Cells.Item(59, 8).Value = 17857.143  # H59 was 5000
$ws.Cells.Item(59, 10).Value = 20000  # J59 was 0
$ws.Cells.Item(59, 12).Value = 20000  # L59 was 0
$ws.Cells.Item(59, 14).Value = -22290  # N59 was None
$ws.Cells.Item(97, 8).Value = 25915.166  # H97 was 27272.75
$ws.Cells.Item(97, 10).Value = 25915.166  # J97 was 27272.75
$ws.Cells.Item(97, 12).Value = 25915.166  # L97 was 27272.75
$ws.Cells.Item(97, 14).Value = -27897.166  # N97 was -29254.75
$ws.Cells.Item(122, 8).Value = 5737.375  # H122 was 5649.875
$ws.Cells.Item(122, 9).Value = 4679.8  # I122 was 5349.75
$ws.Cells.Item(122, 10).Value = 7500  # J122 was 5950
$ws.Cells.Item(122, 11).Value = 14039.4  # K122 was 16049.25
$ws.Cells.Item(122, 12).Value = 22500  # L122 was 17850
$ws.Cells.Item(122, 13).Value = -11589.4  # M122 was -13599.25
$ws.Cells.Item(122, 14).Value = -27400  # N122 was -22750
$ws.Cells.Item(136, 8).Value = 2008.5714  # H136 was 1443.5
$ws.Cells.Item(136, 9).Value = 1510  # I136 was 1443.5
$ws.Cells.Item(136, 10).Value = 5000  # J136 was 0
$ws.Cells.Item(136, 11).Value = 4530  # K136 was 4330.5
$ws.Cells.Item(136, 12).Value = 15000  # L136 was 0
$ws.Cells.Item(136, 13).Value = -1980  # M136 was -1780.5
$ws.Cells.Item(136, 14).Value = -20100  # N136 was None

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 602.25  # H113 was 649.125
$ws.Cells.Item(113, 9).Value = 536.5  # I113 was 560
$ws.Cells.Item(113, 10).Value = 799.5  # J113 was 797.6667
$ws.Cells.Item(113, 11).Value = 1609.5  # K113 was 1680
$ws.Cells.Item(113, 12).Value = 2398.5  # L113 was 2393.0001
$ws.Cells.Item(113, 13).Value = 560.5  # M113 was 490
$ws.Cells.Item(113, 14).Value = -6738.5  # N113 was -6733.0001
$ws.Cells.Item(131, 8).Value = 1589.3448  # H131 was 1561.5
$ws.Cells.Item(131, 10).Value = 2319.5386  # J131 was 2064.111
$ws.Cells.Item(131, 12).Value = 6958.6158  # L131 was 6192.333
$ws.Cells.Item(131, 14).Value = -17038.6158  # N131 was -16272.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 635.2105  # H2 was 583.7
$ws.Cells.Item(2, 9).Value = 289.0909  # I2 was 264.58334
$ws.Cells.Item(2, 10).Value = 1111.125  # J2 was 1062.375
$ws.Cells.Item(2, 11).Value = 289.0909  # K2 was 264.58334
$ws.Cells.Item(2, 12).Value = 1111.125  # L2 was 1062.375
$ws.Cells.Item(2, 13).Value = -176.0909  # M2 was -151.58334
$ws.Cells.Item(2, 14).Value = -1337.125  # N2 was -1288.375
$ws.Cells.Item(20, 8).Value = 35714.363  # H20 was 37085.2
$ws.Cells.Item(20, 9).Value = 6150  # I20 was 5100
$ws.Cells.Item(20, 10).Value = 42284.223  # J20 was 50793.145
$ws.Cells.Item(20, 11).Value = 6150  # K20 was 5100
$ws.Cells.Item(20, 12).Value = 42284.223  # L20 was 50793.145
$ws.Cells.Item(20, 13).Value = -5905  # M20 was -4855
$ws.Cells.Item(20, 14).Value = -42774.223  # N20 was -51283.145
$ws.Cells.Item(68, 8).Value = 0  # H68 was 28295
$ws.Cells.Item(68, 10).Value = 0  # J68 was 28295
$ws.Cells.Item(68, 12).Value = 0  # L68 was 28295
$ws.Cells.Item(68, 14).ClearContents()  # N68 was -29917
$ws.Cells.Item(69, 8).Value = 0  # H69 was 22201
$ws.Cells.Item(69, 10).Value = 0  # J69 was 22201
$ws.Cells.Item(69, 12).Value = 0  # L69 was 22201
$ws.Cells.Item(69, 14).ClearContents()  # N69 was -23699
$ws.Cells.Item(71, 8).Value = 0  # H71 was 28295
$ws.Cells.Item(71, 10).Value = 0  # J71 was 28295
$ws.Cells.Item(71, 12).Value = 0  # L71 was 84885
$ws.Cells.Item(71, 14).ClearContents()  # N71 was -92997
$ws.Cells.Item(72, 8).Value = 0  # H72 was 22201
$ws.Cells.Item(72, 10).Value = 0  # J72 was 22201
$ws.Cells.Item(72, 12).Value = 0  # L72 was 66603
$ws.Cells.Item(72, 14).ClearContents()  # N72 was -74091
$ws.Cells.Item(74, 8).Value = 0  # H74 was 54999
$ws.Cells.Item(74, 10).Value = 0  # J74 was 54999
$ws.Cells.Item(74, 12).Value = 0  # L74 was 54999
$ws.Cells.Item(74, 14).ClearContents()  # N74 was -56871
$ws.Cells.Item(75, 8).Value = 0  # H75 was 45000
$ws.Cells.Item(75, 10).Value = 0  # J75 was 45000
$ws.Cells.Item(75, 12).Value = 0  # L75 was 45000
$ws.Cells.Item(75, 14).ClearContents()  # N75 was -46748
$ws.Cells.Item(77, 8).Value = 0  # H77 was 54999
$ws.Cells.Item(77, 10).Value = 0  # J77 was 54999
$ws.Cells.Item(77, 12).Value = 0  # L77 was 164997
$ws.Cells.Item(77, 14).ClearContents()  # N77 was -174357
$ws.Cells.Item(78, 8).Value = 0  # H78 was 45000
$ws.Cells.Item(78, 10).Value = 0  # J78 was 45000
$ws.Cells.Item(78, 12).Value = 0  # L78 was 135000
$ws.Cells.Item(78, 14).ClearContents()  # N78 was -143736
$ws.Cells.Item(80, 8).Value = 10783.77  # H80 was 10096.071
$ws.Cells.Item(80, 9).Value = 18316.666  # I80 was 14006.875
$ws.Cells.Item(80, 10).Value = 4327  # J80 was 4881.6665
$ws.Cells.Item(80, 11).Value = 18316.666  # K80 was 14006.875
$ws.Cells.Item(80, 12).Value = 4327  # L80 was 4881.6665
$ws.Cells.Item(80, 13).Value = -17318.666  # M80 was -13008.875
$ws.Cells.Item(80, 14).Value = -6323  # N80 was -6877.6665
$ws.Cells.Item(83, 8).Value = 10783.77  # H83 was 10096.071
$ws.Cells.Item(83, 9).Value = 18316.666  # I83 was 14006.875
$ws.Cells.Item(83, 10).Value = 4327  # J83 was 4881.6665
$ws.Cells.Item(83, 11).Value = 91583.33  # K83 was 70034.375
$ws.Cells.Item(83, 12).Value = 21635  # L83 was 24408.3325
$ws.Cells.Item(83, 13).Value = -86591.33  # M83 was -65042.375
$ws.Cells.Item(83, 14).Value = -31619  # N83 was -34392.3325

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 409.53333  # H55 was 428.96667
$ws.Cells.Item(55, 9).Value = 311.73334  # I55 was 352.64285
$ws.Cells.Item(55, 10).Value = 507.33334  # J55 was 495.75
$ws.Cells.Item(55, 11).Value = 311.73334  # K55 was 352.64285
$ws.Cells.Item(55, 12).Value = 507.33334  # L55 was 495.75
$ws.Cells.Item(55, 13).Value = -138.73334  # M55 was -179.64285
$ws.Cells.Item(55, 14).Value = -853.33334  # N55 was -841.75
$ws.Cells.Item(122, 8).Value = 15996.7  # H122 was 10976.667
$ws.Cells.Item(122, 9).Value = 18497.125  # I122 was 11743.077
$ws.Cells.Item(122, 11).Value = 55491.375  # K122 was 35229.231
$ws.Cells.Item(122, 13).Value = -53041.375  # M122 was -32779.231
$ws.Cells.Item(132, 8).Value = 2975.0454  # H132 was 2849.6
$ws.Cells.Item(132, 9).Value = 2576.6316  # I132 was 2464.2856
$ws.Cells.Item(132, 10).Value = 5498.3335  # J132 was 4872.5
$ws.Cells.Item(132, 11).Value = 7729.8948  # K132 was 7392.8568
$ws.Cells.Item(132, 12).Value = 16495.0005  # L132 was 14617.5
$ws.Cells.Item(132, 13).Value = -5199.8948  # M132 was -4862.8568
$ws.Cells.Item(132, 14).Value = -21555.0005  # N132 was -19677.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 1185325.6  # H14 was 1422300.8
$ws.Cells.Item(14, 10).Value = 450  # J14 was 0
$ws.Cells.Item(14, 12).Value = 450  # L14 was 0
$ws.Cells.Item(14, 14).Value = -786  # N14 was None
$ws.Cells.Item(136, 8).Value = 1989  # H136 was 2056.75
$ws.Cells.Item(136, 9).Value = 1952.8928  # I136 was 2021.8148
$ws.Cells.Item(136, 11).Value = 5858.678400000001  # K136 was 6065.4444
$ws.Cells.Item(136, 13).Value = -3308.678400000001  # M136 was -3515.4444
